$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

$ws_ALC.Range("H9").Value = 118.63636
$ws_ALC.Range("I9").Value = 87.625
$ws_ALC.Range("J9").Value = 201.33333
$ws_ALC.Range("K9").Value = 87.625
$ws_ALC.Range("L9").Value = 201.33333
$ws_ALC.Range("M9").Value = 81.375
$ws_ALC.Range("N9").Value = -539.3333299999999

$ws_ALC.Range("H18").Value = 2967
$ws_ALC.Range("I18").Value = 3400.5
$ws_ALC.Range("J18").Value = 2100
$ws_ALC.Range("K18").Value = 3400.5
$ws_ALC.Range("L18").Value = 2100
$ws_ALC.Range("M18").Value = -3116.5
$ws_ALC.Range("N18").Value = -2668

$ws_ALC.Range("H19").Value = 1358.6
$ws_ALC.Range("J19").Value = 1138
$ws_ALC.Range("L19").Value = 1138
$ws_ALC.Range("N19").Value = -1488

$ws_ALC.Range("H33").Value = 125.4375
$ws_ALC.Range("I33").Value = 110.53333
$ws_ALC.Range("K33").Value = 110.53333
$ws_ALC.Range("M33").Value = 118.46667

$ws_ALC.Range("H40").Value = 2127.182
$ws_ALC.Range("I40").Value = 1377.6666
$ws_ALC.Range("K40").Value = 1377.6666
$ws_ALC.Range("M40").Value = -1202.6666

$ws_ALC.Range("H70").Value = 1959.3334
$ws_ALC.Range("I70").Value = 1798.875
$ws_ALC.Range("K70").Value = 5396.625
$ws_ALC.Range("M70").Value = -5126.625

$ws_ALC.Range("H73").Value = 1959.3334
$ws_ALC.Range("I73").Value = 1798.875
$ws_ALC.Range("K73").Value = 5396.625
$ws_ALC.Range("M73").Value = -4460.625

$ws_ALC.Range("H80").Value = 1151.875
$ws_ALC.Range("I80").Value = 788
$ws_ALC.Range("J80").Value = 1758.3334
$ws_ALC.Range("K80").Value = 2364
$ws_ALC.Range("L80").Value = 5275.0002
$ws_ALC.Range("M80").Value = -1366
$ws_ALC.Range("N80").Value = -7271.0002

$ws_ALC.Range("H83").Value = 1151.875
$ws_ALC.Range("I83").Value = 788
$ws_ALC.Range("J83").Value = 1758.3334
$ws_ALC.Range("K83").Value = 7092
$ws_ALC.Range("L83").Value = 15825.0006
$ws_ALC.Range("M83").Value = -2100
$ws_ALC.Range("N83").Value = -25809.0006

$ws_ALC.Range("H92").Value = 363.22223
$ws_ALC.Range("I92").Value = 327.375
$ws_ALC.Range("K92").Value = 327.375
$ws_ALC.Range("M92").Value = 920.625

$ws_ALC.Range("H100").Value = 4011.1333
$ws_ALC.Range("J100").Value = 3033.3333
$ws_ALC.Range("L100").Value = 3033.3333
$ws_ALC.Range("N100").Value = -4115.3333

$ws_ALC.Range("H138").Value = 3110.25
$ws_ALC.Range("I138").Value = 1470.5
$ws_ALC.Range("K138").Value = 4411.5
$ws_ALC.Range("M138").Value = 728.5

$ws_ARM.Range("H92").Value = 50000
$ws_ARM.Range("J92").Value = 50000
$ws_ARM.Range("L92").Value = 50000
$ws_ARM.Range("N92").Value = -54992

$ws_ARM.Range("H122").Value = 2711.3333
$ws_ARM.Range("I122").Value = 2700.2856
$ws_ARM.Range("K122").Value = 8100.8568
$ws_ARM.Range("M122").Value = -5650.8568

$ws_ARM.Range("H134").Value = 85000
$ws_ARM.Range("J134").Value = 85000
$ws_ARM.Range("L134").Value = 85000
$ws_ARM.Range("N134").Value = -95140

$ws_BSM.Range("H20").Value = 4999.5
$ws_BSM.Range("I20").Value = 4999.5
$ws_BSM.Range("K20").Value = 4999.5
$ws_BSM.Range("M20").Value = -4752.5

$ws_BSM.Range("H61").Value = 50001
$ws_BSM.Range("J61").Value = 50001
$ws_BSM.Range("L61").Value = 50001
$ws_BSM.Range("N61").Value = -50627

$ws_BSM.Range("H94").Value = 1679
$ws_BSM.Range("I94").Value = 1461.6923
$ws_BSM.Range("J94").Value = 2385.25
$ws_BSM.Range("K94").Value = 1461.6923
$ws_BSM.Range("L94").Value = 2385.25
$ws_BSM.Range("M94").Value = -1010.6923
$ws_BSM.Range("N94").Value = -3287.25

$ws_BSM.Range("H97").Value = 16499.666
$ws_BSM.Range("I97").Value = 16499.666
$ws_BSM.Range("K97").Value = 16499.666
$ws_BSM.Range("M97").Value = -15508.666

$ws_BSM.Range("H103").Value = 23850
$ws_BSM.Range("J103").Value = 23850
$ws_BSM.Range("L103").Value = 23850
$ws_BSM.Range("N103").Value = -26194

$ws_BSM.Range("H105").Value = 3213.4285
$ws_BSM.Range("I105").Value = 3240
$ws_BSM.Range("J105").Value = 3147
$ws_BSM.Range("K105").Value = 3240
$ws_BSM.Range("L105").Value = 3147
$ws_BSM.Range("M105").Value = -1493
$ws_BSM.Range("N105").Value = -6641

$ws_CRP.Range("H16").Value = 1749
$ws_CRP.Range("I16").Value = 1749
$ws_CRP.Range("K16").Value = 1749
$ws_CRP.Range("M16").Value = -1462

$ws_CRP.Range("H31").Value = 1507.6666
$ws_CRP.Range("I31").Value = 1311.75
$ws_CRP.Range("J31").Value = 1899.5
$ws_CRP.Range("K31").Value = 1311.75
$ws_CRP.Range("L31").Value = 1899.5
$ws_CRP.Range("M31").Value = -1016.75
$ws_CRP.Range("N31").Value = -2489.5

$ws_CRP.Range("H34").Value = 1507.6666
$ws_CRP.Range("I34").Value = 1311.75
$ws_CRP.Range("J34").Value = 1899.5
$ws_CRP.Range("K34").Value = 1311.75
$ws_CRP.Range("L34").Value = 1899.5
$ws_CRP.Range("M34").Value = -1109.75
$ws_CRP.Range("N34").Value = -2303.5

$ws_CRP.Range("H58").Value = 1733.4762
$ws_CRP.Range("J58").Value = 1666.5
$ws_CRP.Range("L58").Value = 1666.5
$ws_CRP.Range("N58").Value = -2072.5

$ws_CRP.Range("H113").Value = 1749
$ws_CRP.Range("I113").Value = 1749
$ws_CRP.Range("K113").Value = 1749
$ws_CRP.Range("M113").Value = 421

$ws_CRP.Range("H132").Value = 3380.2727
$ws_CRP.Range("J132").Value = 3072.2
$ws_CRP.Range("L132").Value = 9216.599999999999
$ws_CRP.Range("N132").Value = -14276.6

$ws_CRP.Range("H136").Value = 1733.4762
$ws_CRP.Range("J136").Value = 1666.5
$ws_CRP.Range("L136").Value = 4999.5
$ws_CRP.Range("N136").Value = -10099.5

$ws_GSM.Range("H58").Value = 9525
$ws_GSM.Range("I58").Value = 9525
$ws_GSM.Range("K58").Value = 9525
$ws_GSM.Range("M58").Value = -9248

$ws_GSM.Range("H113").Value = 105
$ws_GSM.Range("J113").Value = 0
$ws_GSM.Range("L113").Value = 0
$ws_GSM.Range("N113").ClearContents()

$ws_GSM.Range("H126").Value = 9266.666999999999
$ws_GSM.Range("I126").Value = 11000
$ws_GSM.Range("J126").Value = 8400
$ws_GSM.Range("K126").Value = 33000
$ws_GSM.Range("L126").Value = 25200
$ws_GSM.Range("M126").Value = -30530
$ws_GSM.Range("N126").Value = -30140

$ws_GSM.Range("H138").Value = 49673.285
$ws_GSM.Range("J138").Value = 49673.285
$ws_GSM.Range("L138").Value = 49673.285
$ws_GSM.Range("N138").Value = -59953.285

$ws_LTW.Range("H46").Value = 4772.952
$ws_LTW.Range("J46").Value = 3135.7144
$ws_LTW.Range("L46").Value = 3135.7144
$ws_LTW.Range("N46").Value = -3511.7144

$ws_LTW.Range("H132").Value = 2400
$ws_LTW.Range("I132").Value = 2400
$ws_LTW.Range("J132").Value = 0
$ws_LTW.Range("K132").Value = 7200
$ws_LTW.Range("L132").Value = 0
$ws_LTW.Range("M132").Value = -4670
$ws_LTW.Range("N132").ClearContents()

$ws_LTW.Range("H136").Value = 3162.375
$ws_LTW.Range("I136").Value = 2410.3
$ws_LTW.Range("J136").Value = 4415.8335
$ws_LTW.Range("K136").Value = 7230.900000000001
$ws_LTW.Range("L136").Value = 13247.5005
$ws_LTW.Range("M136").Value = -4680.900000000001
$ws_LTW.Range("N136").Value = -18347.5005
